$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.238.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.257.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.255.08"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.785.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.251.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.268.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("E32").Value = "  -0.11%  "

$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.82%  "

$ws.Range("E35").Value = "  -0.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0721"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0393"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "420.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.002.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.38%  "

$ws.Range("E44").Value = "  -7.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.12%  "
